$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark row 12 (intern_.fun_args_pos.R) as done across columns C, E, F, G, H, I, J, K, L, M
$cols = @("C", "E", "F", "G", "H", "I", "J", "K", "L", "M")
foreach ($col in $cols) {
    $addr = $col + "12"
    $ws.Range($addr).Value = "x"
}

# Update the active selection cell shown in the sheet view
$ws.Range("N16").Select()
